{"js": "// Word adds the new \"(${hari_penginapan} hari x Rp ${biaya_penginapan},-)\"\n// qualifier onto the \"Biaya Penginapan\" row label, relocates the document's\n// \"_GoBack\" bookmark into that new text (right after \"biaya_penginapan\"),\n// and renames the placeholder in the amount cell from \"${biaya_penginapan}\"\n// to \"${biaya_penginapan_total}\".\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The \"Biaya Penginapan\" row lives in the first table on the page\n// (row index 4, 0-based): col 1 is the label, col 3 is the amount.\nconst table = tables.items[0];\nconst labelCell = table.getCell(4, 1);\nconst amountCell = table.getCell(4, 3);\n\nconst labelPara = labelCell.body.paragraphs.getFirst();\nconst amountPara = amountCell.body.paragraphs.getFirst();\n\n// 1) Append the day-count / rate breakdown to the row label.\nlabelPara.insertText(\n  \" (${hari_penginapan} hari x Rp ${biaya_penginapan},-)\",\n  \"End\"\n);\nawait context.sync();\n\n// 2) Update the amount cell's placeholder: biaya_penginapan -> biaya_penginapan_total.\nconst amountMatches = amountPara.search(\"${biaya_penginapan}\", { matchCase: true });\namountMatches.load(\"items\");\nawait context.sync();\namountMatches.items[0].insertText(\"${biaya_penginapan_total}\", \"Replace\");\nawait context.sync();\n\n// 3) Move the \"_GoBack\" bookmark from the end of the document into the\n// label cell, right after the newly-inserted \"biaya_penginapan\" token.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst labelMatches = labelPara.search(\"biaya_penginapan\", { matchCase: true });\nlabelMatches.load(\"items\");\nawait context.sync();\nconst bookmarkPoint = labelMatches.items[0].getRange(\"End\");\nbookmarkPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Word adds the new \"(${hari_penginapan} hari x Rp ${biaya_penginapan},-)\"\n# qualifier onto the \"Biaya Penginapan\" row label, relocates the document's\n# \"_GoBack\" bookmark into that new text (right after \"biaya_penginapan\"),\n# and renames the placeholder in the amount cell from \"${biaya_penginapan}\"\n# to \"${biaya_penginapan_total}\".\n\n$d = $word.ActiveDocument\n\n# The \"Biaya Penginapan\" row lives in the first table on the page\n# (row 5, 1-based): col 2 is the label, col 4 is the amount.\n$table = $d.Tables.Item(1)\n$labelCell = $table.Cell(5, 2)\n$amountCell = $table.Cell(5, 4)\n\n# 1) Append the day-count / rate breakdown to the row label.\n$labelRange = $labelCell.Range\n$labelRange.End = $labelRange.End - 1\n$labelRange.InsertAfter(\" (`${hari_penginapan} hari x Rp `${biaya_penginapan},-)\")\n\n# 2) Update the amount cell's placeholder: biaya_penginapan -> biaya_penginapan_total.\n# (Re-fetch a fresh Range via start/end offsets -- a Range returned from\n# .Duplicate can keep stale positions once earlier edits shift the story.)\n$amountCell2 = $d.Tables.Item(1).Cell(5, 4)\n$amountRange = $d.Range($amountCell2.Range.Start, $amountCell2.Range.End)\n$amountRange.Find.Execute(\"biaya_penginapan\") | Out-Null\n$amountPoint = $d.Range($amountRange.End, $amountRange.End)\n$amountPoint.InsertAfter(\"_total\")\n\n# 3) Move the \"_GoBack\" bookmark from the end of the document into the\n# label cell, right after the newly-inserted \"biaya_penginapan\" token.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$labelCell2 = $d.Tables.Item(1).Cell(5, 2)\n$labelFindRange = $d.Range($labelCell2.Range.Start, $labelCell2.Range.End)\n$labelFindRange.Find.Execute(\"biaya_penginapan\") | Out-Null\n$bookmarkPoint = $d.Range($labelFindRange.End, $labelFindRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkPoint)\n"}
